# Rebrand accent color 1D85AD -> 0D3E51 across the remaining not-yet-updated
# shapes, and drop the white (bg1) line/border on shapes whose card has
# already migrated to the darker navy background.
#
# RGB(0x0D,0x3E,0x51) = 5324301   (new navy)
# RGB(0x1D,0x85,0xAD) = 11371805  (old teal)  -- for reference only

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 : "Group 16" card (logos-only band, no title textbox) plus the
# "Group 1" card's title textbox.
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Group 1 -> TextBox 8 ("RESPIMATIC 100 WEB RECORDER") : fill color only.
$group1 = $s1.Shapes.Item("Group 1")
$textBox8 = $group1.GroupItems.Item("TextBox 8")
$textBox8.Fill.ForeColor.RGB = 5324301

# Group 16 -> Rectangle 10 (card background) : fill color + remove white line.
$group16 = $s1.Shapes.Item("Group 16")
$rectangle10 = $group16.GroupItems.Item("Rectangle 10")
$rectangle10.Fill.ForeColor.RGB = 5324301
$rectangle10.Line.Visible = $false

# Group 16 -> the three logo pictures: remove their white line.
$picture11 = $group16.GroupItems.Item("Picture 11")
$picture11.Line.Visible = $false

$picture12 = $group16.GroupItems.Item("Picture 12")
$picture12.Line.Visible = $false

$picture13 = $group16.GroupItems.Item("Picture 13")
$picture13.Line.Visible = $false

# ---------------------------------------------------------------------
# Slide 2 : the two left-column textboxes ("WEB APPS" and the
# "O2 Flow Rate Calculator" box).
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# TextBox 1 ("WEB APPS ... RESPIMATIC 100") : move up slightly, recolor,
# and remove its white line.
$textBox1 = $s2.Shapes.Item("TextBox 1")
$textBox1.Top = 255602 / 914400 * 72
$textBox1.Fill.ForeColor.RGB = 5324301
$textBox1.Line.Visible = $false

# TextBox 4 ("O2 Flow Rate Calculator ... RESPIMATIC 100") : fill color only,
# the white line stays.
$textBox4 = $s2.Shapes.Item("TextBox 4")
$textBox4.Fill.ForeColor.RGB = 5324301
